# The three species-observation records that occupied rows 12-14 were
# re-ordered (the record that used to be first, A=130826137, now sorts
# last). Concretely:
#   new row 12  <=  old row 13
#   new row 13  <=  old row 14
#   new row 14  <=  old row 12
# Every column (A:AY) moves along with its row, so instead of editing
# individual cells we snapshot the three full rows first and then write
# the snapshots back in rotated order. This naturally also removes cells
# that shouldn't exist in their new position (e.g. M/AC on the new row 14)
# and creates cells that are new in their new position (e.g. M/AC on the
# new row 12), because writing a $null value to a Range clears the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 12
$lastRow = 14
$lastCol = 51   # column AY

# Snapshot the rows before touching anything.
$snapshots = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $snapshots += , $rowVals
}

# New order: old row13, old row14, old row12 (indices 1, 2, 0 of $snapshots)
$order = @(1, 2, 0)

# Only touch a cell when its content actually needs to change - e.g. many
# cells (dates, county/commune names, booleans...) happen to already hold
# the correct value for their new row, and re-assigning an unchanged
# date-like text such as "2026-01-22" would make Excel's Value setter
# "helpfully" re-interpret it as a real date serial, which is not what the
# source workbook has (those cells are plain inline strings there).
for ($i = 0; $i -lt 3; $i++) {
    $destRow = $firstRow + $i
    $srcVals = $snapshots[$order[$i]]
    for ($c = 1; $c -le $lastCol; $c++) {
        $newVal = $srcVals[$c - 1]
        $curVal = $ws.Cells.Item($destRow, $c).Value()
        if ($curVal -ne $newVal) {
            $ws.Cells.Item($destRow, $c).Value = $newVal
        }
    }
}
